# Apply the "daily refresh" update to the Saldo_guide workbook:
#  - rename the sheet to the new run's timestamped name
#  - bump the reference date (column G) for every data row by one day
#  - update the handful of balances (columns E and H) that changed between runs

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename worksheet tab to reflect the new export timestamp
$ws.Name = "IClientBalance-20241205-102334-"

# 2) Bump the "Dt. Referencia" date (column G) from 45630 to 45631 for rows 2-274
$lastRow = 274
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 7).Value = 45631
}

# 3) Update the balances that changed for specific accounts (columns E and H)
$changedRows = @{
    6   = 15391.31
    105 = 6284.22
    138 = 9701.11
    143 = 42991.54
    264 = 341.4
}

foreach ($r in $changedRows.Keys) {
    $newVal = $changedRows[$r]
    $ws.Cells.Item($r, 5).Value = $newVal   # column E - Saldo Previsto
    $ws.Cells.Item($r, 8).Value = $newVal   # column H - Vl. Total
}
